$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update shared-string-backed Config Name values for rows 2-4
$ws.Range("A2").Value = "Debug"
$ws.Range("A3").Value = "Test"
$ws.Range("A4").Value = "Fast"

# Row 3 values update (B3..J3) - F3 stays 0 (unchanged), others updated
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 20
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = 20
$ws.Range("G3").Value = 20
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 53
$ws.Range("J3").Value = 3

# Row 4 values update (B4..J4)
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 8
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = 8
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 50
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 45
$ws.Range("J4").Value = 3

# F4 lost its cell style (quote-prefix format) when the value was set; restore
# it by copying the number format from F2, which carries the same style.
$ws.Range("F2").Copy()
$ws.Range("F4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# View changes: scroll back to A1 (clear topLeftCell) and select F7
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("F7").Select()
